# Update workbook per commit: "Mise à jour du fichier via Shiny"
# Data refresh across pro / ind / VA / conso sheets, plus an active-sheet/selection change.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("pro")
$ws2 = $wb.Worksheets.Item("ind")
$ws3 = $wb.Worksheets.Item("VA")
$ws4 = $wb.Worksheets.Item("conso")

# --- Sheet "pro" (sheet1): refresh annual production values B2:B26 ---
$ws1.Range("B2").Value = 2177088.076454516
$ws1.Range("B3").Value = 2567810.8130529216
$ws1.Range("B4").Value = 1988303.522599349
$ws1.Range("B5").Value = 2469849.8921814561
$ws1.Range("B6").Value = 2473548.5959649342
$ws1.Range("B7").Value = 2348682.1188426386
$ws1.Range("B8").Value = 1880535.4255204957
$ws1.Range("B9").Value = 1603967.9752063858
$ws1.Range("B10").Value = 1509900.0629585646
$ws1.Range("B11").Value = 1593852.274622286
$ws1.Range("B12").Value = 1539984.947953685
$ws1.Range("B13").Value = 1627491.0874458575
$ws1.Range("B14").Value = 1855313.5995463191
$ws1.Range("B15").Value = 1968075.5539721963
$ws1.Range("B16").Value = 1937562.8742623129
$ws1.Range("B17").Value = 1447925.5507954138
$ws1.Range("B18").Value = 1713408.5714997195
$ws1.Range("B19").Value = 1872780.1695488854
$ws1.Range("B20").Value = 2024197.3832623838
$ws1.Range("B21").Value = 2287047
$ws1.Range("B22").Value = 2336188
$ws1.Range("B23").Value = 2415516.6566253239
$ws1.Range("B24").Value = 2526122.0189849623
$ws1.Range("B25").Value = 2632110.7650183784
$ws1.Range("B26").Value = 2574511.2535821563

# --- Sheet "ind" (sheet2): refresh quarterly/monthly index values B2:B101 ---
$ws2.Range("B2").Value = 413489.19559666933
$ws2.Range("B3").Value = 470400.84892142564
$ws2.Range("B4").Value = 778201.61142732284
$ws2.Range("B5").Value = 701851.96014975791
$ws2.Range("B6").Value = 623494.21256565175
$ws2.Range("B7").Value = 626099.11946481769
$ws2.Range("B8").Value = 810539.52292616444
$ws2.Range("B9").Value = 728068.52490767534
$ws2.Range("B10").Value = 507972.36413364322
$ws2.Range("B11").Value = 482500.26872867747
$ws2.Range("B12").Value = 607224.08571075113
$ws2.Range("B13").Value = 561259.30619319517
$ws2.Range("B14").Value = 616584.76314718928
$ws2.Range("B15").Value = 618781.12605550338
$ws2.Range("B16").Value = 764760.81792595191
$ws2.Range("B17").Value = 681705.94285954721
$ws2.Range("B18").Value = 676418.50144925714
$ws2.Range("B19").Value = 691581.09604557208
$ws2.Range("B20").Value = 753034.95194456703
$ws2.Range("B21").Value = 564814.25740131037
$ws2.Range("B22").Value = 538453.10882014642
$ws2.Range("B23").Value = 670479.10522003472
$ws2.Range("B24").Value = 826483.2060144376
$ws2.Range("B25").Value = 514849.84533458378
$ws2.Range("B26").Value = 519876.68992382934
$ws2.Range("B27").Value = 526642.27077573293
$ws2.Range("B28").Value = 632224.13109540427
$ws2.Range("B29").Value = 363195.29465729825
$ws2.Range("B30").Value = 314643.58468187408
$ws2.Range("B31").Value = 402862.12370343308
$ws2.Range("B32").Value = 532356.62815303018
$ws2.Range("B33").Value = 491771.31458393292
$ws2.Range("B34").Value = 406831.70267532597
$ws2.Range("B35").Value = 381200.98001170921
$ws2.Range("B36").Value = 520186.17439868656
$ws2.Range("B37").Value = 331273.20268679509
$ws2.Range("B38").Value = 384904.68185258372
$ws2.Range("B39").Value = 427205.16331279383
$ws2.Range("B40").Value = 497327.65151074994
$ws2.Range("B41").Value = 421212.24154122442
$ws2.Range("B42").Value = 347291.68171934405
$ws2.Range("B43").Value = 399389.90411012783
$ws2.Range("B44").Value = 507628.2900977722
$ws2.Range("B45").Value = 417849.20021192008
$ws2.Range("B46").Value = 384247.54748803633
$ws2.Range("B47").Value = 402194.81202630652
$ws2.Range("B48").Value = 507061.01735105185
$ws2.Range("B49").Value = 473672.33275721344
$ws2.Range("B50").Value = 465077.31944581133
$ws2.Range("B51").Value = 473666.37230541476
$ws2.Range("B52").Value = 577465.193130792
$ws2.Range("B53").Value = 498342.93180177972
$ws2.Range("B54").Value = 446894.89742826147
$ws2.Range("B55").Value = 483912.34858680883
$ws2.Range("B56").Value = 613871.09522048361
$ws2.Range("B57").Value = 592313.58487072482
$ws2.Range("B58").Value = 478843.82801656512
$ws2.Range("B59").Value = 528823.82574548945
$ws2.Range("B60").Value = 660406.27399803768
$ws2.Range("B61").Value = 435786.47041663935
$ws2.Range("B62").Value = 244012.88885768989
$ws2.Range("B63").Value = 289627.66286910599
$ws2.Range("B64").Value = 496774.34451129287
$ws2.Range("B65").Value = 541783.4908653457
$ws2.Range("B66").Value = 463261.11223519279
$ws2.Range("B67").Value = 514557.7748138727
$ws2.Range("B68").Value = 652932.6831094576
$ws2.Range("B69").Value = 634240.8113208327
$ws2.Range("B70").Value = 549904.40717999113
$ws2.Range("B71").Value = 633574.62089773279
$ws2.Range("B72").Value = 825855.89640463633
$ws2.Range("B73").Value = 803111.82501765271
$ws2.Range("B74").Value = 712515.43318206398
$ws2.Range("B75").Value = 779712.57461952139
$ws2.Range("B76").Value = 895051.60911139287
$ws2.Range("B77").Value = 855922.55673974112
$ws2.Range("B78").Value = 861136.4602990424
$ws2.Range("B79").Value = 974229.16675269336
$ws2.Range("B80").Value = 1150219.0582765213
$ws2.Range("B81").Value = 1014415.3146717428
$ws2.Range("B82").Value = 979889.03249605314
$ws2.Range("B83").Value = 1021556.2153226548
$ws2.Range("B84").Value = 1331554.9903030754
$ws2.Range("B85").Value = 1238068.1028300705
$ws2.Range("B86").Value = 1105474.1858093182
$ws2.Range("B87").Value = 1176839.1466706307
$ws2.Range("B88").Value = 1499209.2429791924
$ws2.Range("B89").Value = 1314933.2639614381
$ws2.Range("B90").Value = 1165783.6890707607
$ws2.Range("B91").Value = 1233999.3861137838
$ws2.Range("B92").Value = 1573552.7907735379
$ws2.Range("B93").Value = 1390366.9415390373
$ws2.Range("B94").Value = 1254871.5174689167
$ws2.Range("B95").Value = 1288252.6128015923
$ws2.Range("B96").Value = 1613970.1783435808
$ws2.Range("B97").Value = 1439466.8274332813
$ws2.Range("B98").Value = 1058543.5480696305
$ws2.Range("B100").Value = 479904.86140073417
$ws2.Range("B101").Value = 719325.99147441436

# Row 99 becomes a formula referencing the (now-updated) B100
$ws2.Range("B99").Formula = "=B100"

# --- Sheet "conso" (sheet4): refresh annual consumption values B2:B26 ---
$ws4.Range("B2").Value = 1571477.7391109124
$ws4.Range("B3").Value = 1853486.9052523435
$ws4.Range("B4").Value = 1435191.9058242077
$ws4.Range("B5").Value = 1782792.3846726893
$ws4.Range("B6").Value = 1785466.3045179155
$ws4.Range("B7").Value = 1695323.8162250468
$ws4.Range("B8").Value = 1357397.2246882429
$ws4.Range("B9").Value = 1157755.8038547977
$ws4.Range("B10").Value = 1089857.4754876068
$ws4.Range("B11").Value = 1150453.1112369266
$ws4.Range("B12").Value = 1111567.1106946126
$ws4.Range("B13").Value = 1174736.6049298833
$ws4.Range("B14").Value = 1339179.0348299067
$ws4.Range("B15").Value = 1420560.5302648027
$ws4.Range("B16").Value = 1398539.8850631886
$ws4.Range("B17").Value = 1045118.5674655235
$ws4.Range("B18").Value = 1236746.9467946021
$ws4.Range("B19").Value = 1351777.3121834856
$ws4.Range("B20").Value = 1461069.7975799066
$ws4.Range("B21").Value = 1650793
$ws4.Range("B22").Value = 1691949.9999999998
$ws4.Range("B23").Value = 1742942.7357080011
$ws4.Range("B24").Value = 1822751.1768937916
$ws4.Range("B25").Value = 1899228.5244320289
$ws4.Range("B26").Value = 1857666.962294556

# "VA" (sheet3) cells are formulas (=pro!Bn -conso!Bn) and recalc automatically.

# --- View / selection state ---
# "pro" loses tabSelected and its lingering selection moves to E88
$ws1.Range("E88").Select()

# "VA" and "conso" selections also move to E88
$ws3.Range("E88").Select()
$ws4.Range("E88").Select()

# "ind" becomes the active sheet, scrolled further down, with D101 selected
$ws2.Activate()
$ws2.Range("D101").Select()
